$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("It decays, leading to entropy and collapse. ", $true, $false, $false, $false, $false, $true, 1, $false, "It decays, leading to entropy and collapse…", 2)
$insStart = $rng.End
$ins = $d.Range($insStart, $insStart)
$ins.Text = "the Shadow Dimensions."
$newRng = $d.Range($insStart, $insStart + 23)
$newRng.Font.Bold = $true
$newRng.Font.NameAscii = "Roboto"
$newRng.Font.NameFarEast = "Roboto"
$newRng.Font.NameOther = "Roboto"
$newRng.Font.NameBi = "Roboto"

$s = $d.Styles.Item(2)
$s.NameLocal = "TableNormal"
